# Add a "Run 50" results column to the PSO table, inserted right before the
# existing "Mean" column, and recompute the "Mean" column so it now averages
# Run 0 .. Run 50 (51 samples) instead of Run 0 .. Run 49 (50 samples).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Mean" column currently lives in column AZ (52nd column). Insert a new,
# blank column there: this pushes the old "Mean" column (with its header,
# values and styling) one position to the right, into column BA, and leaves
# a fresh column AZ for the new "Run 50" data.
$ws.Columns("AZ:AZ").Insert()

# New run data value (constant across every MaxFES row, same as the other
# "Run N" columns in this table) and the corresponding recalculated mean.
$runValue = 2647022.30120277
$meanValue = 1262655.34012531

# Header for the newly inserted column.
$ws.Range("AZ1").Value = "Run 50"

# Fill in the "Run 50" value and the updated "Mean" value for every data row.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = $runValue
    $ws.Cells.Item($r, 53).Value = $meanValue
}
